$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old percent-change/direction columns (I:L) and shift DiD/tstat/pval/Significant (M:P) left into I:L
$ws.Range("I1:L8").Delete()

# Row 2: LPE
$ws.Range("A2").Value = 'LPE'
$ws.Range("B2").Value = 0.7317460317460318
$ws.Range("C2").Value = 0.7744996549344375
$ws.Range("D2").Value = 0.6865079365079364
$ws.Range("E2").Value = 0.7148809523809525
$ws.Range("F2").Value = -0.04523809523809525
$ws.Range("G2").Value = -0.05961870255348516
$ws.Range("H2").Value = 0.01438060731538991
$ws.Range("I2").Value = -0.05726551226551235
$ws.Range("J2").Value = 0.1232615817076855
$ws.Range("K2").Value = 0.904120734233601
$ws.Range("L2").Value = $false

# Row 3: avg_loan_size
$ws.Range("A3").Value = 'avg_loan_size'
$ws.Range("B3").Value = 7415.334802896319
$ws.Range("C3").Value = 7130.230708522579
$ws.Range("D3").Value = 7331.953192307691
$ws.Range("E3").Value = 6983.24620406522
$ws.Range("F3").Value = -83.38161058862816
$ws.Range("G3").Value = -146.9845044573592
$ws.Range("H3").Value = 63.60289386873104
$ws.Range("I3").Value = -136.5767581879309
$ws.Range("J3").Value = 0.05818256419115696
$ws.Range("K3").Value = 0.9548175343646363
$ws.Range("L3").Value = $false

# Row 4: dq29_pot30_payment_rate_unit_per_day
$ws.Range("A4").Value = 'dq29_pot30_payment_rate_unit_per_day'
$ws.Range("B4").Value = 0.008137904442032296
$ws.Range("C4").Value = 0.007757556580780511
$ws.Range("D4").Value = 0.008441831814255781
$ws.Range("E4").Value = 0.007716867004573938
$ws.Range("F4").Value = 0.0003039273722234854
$ws.Range("G4").Value = -0.00004068957620657139
$ws.Range("H4").Value = 0.0003446169484300568
$ws.Range("I4").Value = 0.00001570228808198126
$ws.Range("J4").Value = 0.2884523070102074
$ws.Range("K4").Value = 0.7787000908455818
$ws.Range("L4").Value = $false

# Row 5: dq29_pot30_payment_rate_unit_up_to_day
$ws.Range("A5").Value = 'dq29_pot30_payment_rate_unit_up_to_day'
$ws.Range("B5").Value = 0.7128576856380868
$ws.Range("C5").Value = 0.7104298268035751
$ws.Range("D5").Value = 0.4640292836223756
$ws.Range("E5").Value = 0.478202930382025
$ws.Range("F5").Value = -0.2488284020157112
$ws.Range("G5").Value = -0.2322268964215504
$ws.Range("H5").Value = -0.01660150559416076
$ws.Range("I5").Value = -0.2349435064278677
$ws.Range("J5").Value = -1.171450080479241
$ws.Range("K5").Value = 0.2678242908720487
$ws.Range("L5").Value = $false

# Row 6: dq29_pot30_payment_rate_$_up_to_day
$ws.Range("A6").Value = 'dq29_pot30_payment_rate_$_up_to_day'
$ws.Range("B6").Value = 0.001312028653511341
$ws.Range("C6").Value = -0.00000233055725748332
$ws.Range("D6").Value = 0.004360121972890689
$ws.Range("E6").Value = 0.002091279540988891
$ws.Range("F6").Value = 0.003048093319379348
$ws.Range("G6").Value = 0.002093610098246374
$ws.Range("H6").Value = 0.0009544832211329738
$ws.Range("I6").Value = 0.002249798261704497
$ws.Range("J6").Value = 0.9634707020860324
$ws.Range("K6").Value = 0.355028322995343
$ws.Range("L6").Value = $false

# Row 7: dq30_pct_unit
$ws.Range("A7").Value = 'dq30_pct_unit'
$ws.Range("B7").Value = 0.04669864918067355
$ws.Range("C7").Value = 0.0488541886567222
$ws.Range("D7").Value = 0.03318061298945178
$ws.Range("E7").Value = 0.03262297391715657
$ws.Range("F7").Value = -0.01351803619122178
$ws.Range("G7").Value = -0.01623121473956563
$ws.Range("H7").Value = 0.002713178548343848
$ws.Range("I7").Value = -0.01578724006801846
$ws.Range("J7").Value = 1.448937505851671
$ws.Range("K7").Value = 0.1709773730717109
$ws.Range("L7").Value = $false

# Row 8: dq30_pct_$
$ws.Range("A8").Value = 'dq30_pct_$'
$ws.Range("B8").Value = 0.9981070051443941
$ws.Range("C8").Value = 0.9985858274665315
$ws.Range("D8").Value = 0.9960757709728522
$ws.Range("E8").Value = 0.9978669687437893
$ws.Range("F8").Value = -0.002031234171541952
$ws.Range("G8").Value = -0.0007188587227419011
$ws.Range("H8").Value = -0.001312375448800051
$ws.Range("I8").Value = -0.0009336110689091015
$ws.Range("J8").Value = -1.206522239084651
$ws.Range("K8").Value = 0.2471353352948185
$ws.Range("L8").Value = $false
